# "creating files and section page based on excel spreadsheet"
#
# The Polar Coordinates section had three topic titles re-worded (and the
# "Arclength" row re-labelled/re-ordered), and two worksheets ("Polar
# Coordinates" and "Infinite Series and Sequences") each had a stray,
# completely empty trailing row that is removed. The previously-active
# sheet ("Infinite Series and Sequences") is swapped out in favor of
# "Polar Coordinates" becoming the active/selected tab.

$wb = $excel.ActiveWorkbook

$wsPolar = $wb.Worksheets.Item("Polar Coordinates")
$wsSeries = $wb.Worksheets.Item("Infinite Series and Sequences")

# --- Rename the three polar-curve topics on the "Polar Coordinates" sheet.
# Row 4 = Arclength, Row 5 = Surface Area, Row 6 = Conic Sections (column A).
# Written in (Surface Area, Arc Length, Conic Sections) order so the shared
# string table regenerates with that same ordering.
$wsPolar.Cells.Item(5, 1).Value = "Surface Area with Polar Curves"
$wsPolar.Cells.Item(4, 1).Value = "Arc Length of Polar Curves"
$wsPolar.Cells.Item(6, 1).Value = "Conic Sections with Polar Functions"

# --- Drop the trailing blank row on "Polar Coordinates" (was row 7, a lone
# styled-but-empty B7 cell) so the sheet's used range shrinks to A1:B6.
$wsPolar.Rows.Item(7).Delete()

# --- Drop the trailing blank row on "Infinite Series and Sequences" (was
# row 15, a lone styled-but-empty B15 cell) so the range shrinks to A1:B14.
$wsSeries.Rows.Item(15).Delete()

# --- Switch the active/selected sheet from "Infinite Series and Sequences"
# to "Polar Coordinates", preserving each sheet's own last-known selection.
$wsSeries.Activate()
$wsSeries.Range("B29").Select()

$wsPolar.Activate()
$wsPolar.Range("A9").Select()
